{"js": "// Update the date line and every \"NNN\u00d7N=\" practice-problem cell.\n// Order-sensitive: \"433\u00d72=\" appears twice in the original document (row 1\n// col 1 -> 437\u00d72=, and row 3 col 3 -> 741\u00d74=), so replacements are applied\n// in document order against each search hit rather than a blind find/replace\n// that could collide on duplicate source text.\n\nconst body = context.document.body;\n\n// 1) Date heading.\nconst dateHits = body.search(\"2025-10-09 Thursday\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"2025-10-10 Friday\", \"Replace\");\n}\n\n// 2) Table cells, in the exact document order they appear.\n//    [searchText, [replacement-per-occurrence in document order]]\nconst replacements = [\n  [\"433\u00d72=\", [\"437\u00d72=\", \"741\u00d74=\"]],\n  [\"456\u00d77=\", [\"639\u00d73=\"]],\n  [\"736\u00d79=\", [\"729\u00d78=\"]],\n  [\"495\u00d75=\", [\"626\u00d78=\"]],\n  [\"541\u00d77=\", [\"222\u00d72=\"]],\n  [\"672\u00d73=\", [\"548\u00d72=\"]],\n  [\"753\u00d75=\", [\"597\u00d76=\"]],\n  [\"958\u00d72=\", [\"788\u00d75=\"]],\n  [\"741\u00d77=\", [\"342\u00d79=\"]],\n  [\"782\u00d74=\", [\"772\u00d73=\"]],\n  [\"430\u00d72=\", [\"734\u00d76=\"]],\n  [\"427\u00d74=\", [\"225\u00d73=\"]],\n  [\"257\u00d78=\", [\"621\u00d74=\"]],\n  [\"497\u00d74=\", [\"383\u00d75=\"]],\n  [\"891\u00d79=\", [\"967\u00d76=\"]],\n  [\"616\u00d78=\", [\"306\u00d79=\"]],\n  [\"238\u00d79=\", [\"554\u00d72=\"]],\n  [\"571\u00d76=\", [\"261\u00d76=\"]],\n  [\"182\u00d79=\", [\"984\u00d77=\"]],\n  [\"885\u00d79=\", [\"904\u00d77=\"]],\n  [\"505\u00d76=\", [\"445\u00d76=\"]],\n  [\"707\u00d75=\", [\"511\u00d73=\"]],\n  [\"540\u00d79=\", [\"594\u00d79=\"]],\n  [\"718\u00d79=\", [\"473\u00d72=\"]],\n];\n\nfor (const [needle, repls] of replacements) {\n  const hits = body.search(needle, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < hits.items.length && i < repls.length; i++) {\n    hits.items[i].insertText(repls[i], \"Replace\");\n  }\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and every \"NNN\u00d7N=\" practice-problem cell.\n# Cells are addressed by (row, column) rather than find/replace-by-text\n# because \"433\u00d72=\" appears twice in the source with two different\n# replacement values (row 1 col 1 -> 437\u00d72=, row 10 col 3 -> 741\u00d74=).\n\n$d = $word.ActiveDocument\n\n# 1) Date heading (first paragraph, outside the table) - plain text replace.\n$find = $d.Content.Find\n$find.Text = \"2025-10-09 Thursday\"\n$find.Replacement.Text = \"2025-10-10 Friday\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) Table cells - only the five \"problem\" rows (1, 5, 10, 15, 20) carry\n#    text; the rows in between are blank student-work rows and are left\n#    untouched.\n$t = $d.Tables.Item(1)\n\n$grid = @{\n    1  = @(\"437\u00d72=\", \"639\u00d73=\", \"729\u00d78=\", \"626\u00d78=\", \"222\u00d72=\");\n    5  = @(\"548\u00d72=\", \"597\u00d76=\", \"788\u00d75=\", \"342\u00d79=\", \"772\u00d73=\");\n    10 = @(\"734\u00d76=\", \"225\u00d73=\", \"741\u00d74=\", \"621\u00d74=\", \"383\u00d75=\");\n    15 = @(\"967\u00d76=\", \"306\u00d79=\", \"554\u00d72=\", \"261\u00d76=\", \"984\u00d77=\");\n    20 = @(\"904\u00d77=\", \"445\u00d76=\", \"511\u00d73=\", \"594\u00d79=\", \"473\u00d72=\");\n}\n\nforeach ($r in $grid.Keys) {\n    $values = $grid[$r]\n    for ($c = 1; $c -le $values.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $values[$c - 1]\n    }\n}\n"}
